$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.07488789237668161
$ws.Range("H2").Value2 = 0.07488789237668161
$ws.Range("I2").Value2 = 0.2035874439461883
$ws.Range("J2").Value2 = 0.1610271377764033
$ws.Range("K2").Value2 = 3.67
$ws.Range("L2").Value2 = 0.1645739910313901
$ws.Range("U2").Value2 = 3.4
$ws.Range("V2").Value2 = 0.1538461538461538
$ws.Range("W2").Value2 = 0.2602836879432624
$ws.Range("X2").Value2 = 0.08527046641588971
$ws.Range("Y2").Value2 = 0.1750132215273727
$ws.Range("Z2").Value2 = 1.747922871923499
$ws.Range("AA2").Value2 = 0.2814630171197517
$ws.Range("AB2").Value2 = 0.06874994856208377
$ws.Range("AC2").Value2 = 0.212713068557668
$ws.Range("AD2").Value2 = 18.9
$ws.Range("AF2").Value2 = 18.9
$ws.Range("AG2").Value2 = 15.5
$ws.Range("AH2").Value2 = 0.4609756097560975
$ws.Range("AI2").Value2 = 0.5206611570247934
$ws.Range("AJ2").Value2 = 0.4122340425531914
$ws.Range("AK2").Value2 = 0.4711246200607903
$ws.Range("AL2").Value2 = 0.211
$ws.Range("AM2").Value2 = 0.07099999999999998
$ws.Range("AN2").Value2 = 4.029850746268656
$ws.Range("AO2").Value2 = 21.51658767772512
$ws.Range("AP2").Value2 = 3.304904051172707
$ws.Range("AQ2").Value2 = 63.94366197183101
$ws.Range("G3").Value2 = 0.07488789237668161
$ws.Range("H3").Value2 = 0.07488789237668161
$ws.Range("I3").Value2 = 0.2035874439461883
$ws.Range("J3").Value2 = 0.1610271377764033
$ws.Range("K3").Value2 = 3.67
$ws.Range("L3").Value2 = 0.1645739910313901
$ws.Range("U3").Value2 = 3.4
$ws.Range("V3").Value2 = 0.1538461538461538
$ws.Range("W3").Value2 = 0.2602836879432624
$ws.Range("X3").Value2 = 0.08527046641588971
$ws.Range("Y3").Value2 = 0.1750132215273727
$ws.Range("Z3").Value2 = 1.747922871923499
$ws.Range("AA3").Value2 = 0.2814630171197517
$ws.Range("AB3").Value2 = 0.06874994856208377
$ws.Range("AC3").Value2 = 0.212713068557668
$ws.Range("AD3").Value2 = 18.9
$ws.Range("AF3").Value2 = 18.9
$ws.Range("AG3").Value2 = 15.5
$ws.Range("AH3").Value2 = 0.4609756097560975
$ws.Range("AI3").Value2 = 0.5206611570247934
$ws.Range("AJ3").Value2 = 0.4122340425531914
$ws.Range("AK3").Value2 = 0.4711246200607903
$ws.Range("AL3").Value2 = 0.211
$ws.Range("AM3").Value2 = 0.07099999999999998
$ws.Range("AN3").Value2 = 4.029850746268656
$ws.Range("AO3").Value2 = 21.51658767772512
$ws.Range("AP3").Value2 = 3.304904051172707
$ws.Range("AQ3").Value2 = 63.94366197183101
